# Week 16 stat log + season-sim update (Players Data: Rushing + Receiving)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Rushing sheet
# ---------------------------------------------------------------
$rush = $wb.Worksheets.Item("Rushing")

# M.Rudolph (row 3)
$rush.Cells.Item(3, 4).Value = 1

# N.Harris (row 4)
$rush.Cells.Item(4, 3).Value = 155
$rush.Cells.Item(4, 4).Value = 92
$rush.Cells.Item(4, 5).Value = 23

# B.Snell (row 5)
$rush.Cells.Item(5, 3).Value = 14
$rush.Cells.Item(5, 4).Value = 7
$rush.Cells.Item(5, 5).Value = 2

# ---------------------------------------------------------------
# Receiving sheet
# ---------------------------------------------------------------
$rec = $wb.Worksheets.Item("Receiving")

# N.Harris (row 2)
$rec.Cells.Item(2, 3).Value = 84
$rec.Cells.Item(2, 4).Value = 64
$rec.Cells.Item(2, 7).Value = 14
$rec.Cells.Item(2, 8).Value = 8

# B.Snell (row 3)
$rec.Cells.Item(3, 3).Value = 4
$rec.Cells.Item(3, 4).Value = 2

# D.Johnson (row 7)
$rec.Cells.Item(7, 3).Value = 114
$rec.Cells.Item(7, 4).Value = 79
$rec.Cells.Item(7, 5).Value = 35
$rec.Cells.Item(7, 7).Value = 18
$rec.Cells.Item(7, 8).Value = 9

# C.Claypool (row 8)
$rec.Cells.Item(8, 3).Value = 62
$rec.Cells.Item(8, 4).Value = 39
$rec.Cells.Item(8, 5).Value = 29
$rec.Cells.Item(8, 6).Value = 13
$rec.Cells.Item(8, 7).Value = 10
$rec.Cells.Item(8, 8).Value = 5

# J.Washington (row 9)
$rec.Cells.Item(9, 3).Value = 33

# R.McCloud (row 10)
$rec.Cells.Item(10, 3).Value = 40
$rec.Cells.Item(10, 4).Value = 29
$rec.Cells.Item(10, 5).Value = 6
$rec.Cells.Item(10, 7).Value = 5

# Rows 12-15 got re-sorted (P.Freiermuth dropped behind E.Ebron) -
# row 12 becomes P.Freiermuth's new (much smaller) week-16 stat line,
# row 13 becomes E.Ebron (previously row 12, stats unchanged),
# row 14 (Z.Gentry) now carries P.Freiermuth's former season total,
# row 15 (K.Rader) gets its own updated total.

# Row 12 -> P.Freiermuth
$rec.Cells.Item(12, 2).Value = "P.Freiermuth"
$rec.Cells.Item(12, 3).Value = 3
$rec.Cells.Item(12, 4).Value = 3
$rec.Cells.Item(12, 5).Value = 0
$rec.Cells.Item(12, 6).Value = 0
$rec.Cells.Item(12, 7).Value = 0
$rec.Cells.Item(12, 8).Value = 0

# Row 13 -> E.Ebron
$rec.Cells.Item(13, 2).Value = "E.Ebron"
$rec.Cells.Item(13, 3).Value = 17
$rec.Cells.Item(13, 4).Value = 12
$rec.Cells.Item(13, 5).Value = 0
$rec.Cells.Item(13, 6).Value = 0
$rec.Cells.Item(13, 7).Value = 3
$rec.Cells.Item(13, 8).Value = 1

# Row 14 -> Z.Gentry (name unchanged), new totals
$rec.Cells.Item(14, 3).Value = 58
$rec.Cells.Item(14, 4).Value = 46
$rec.Cells.Item(14, 5).Value = 5
$rec.Cells.Item(14, 6).Value = 2
$rec.Cells.Item(14, 7).Value = 17
$rec.Cells.Item(14, 8).Value = 11

# Row 15 -> K.Rader (name unchanged), new totals
$rec.Cells.Item(15, 3).Value = 15
$rec.Cells.Item(15, 4).Value = 13
$rec.Cells.Item(15, 5).Value = 2
$rec.Cells.Item(15, 6).Value = 1
$rec.Cells.Item(15, 7).Value = 2
$rec.Cells.Item(15, 8).Value = 1

# New row 16 -> C.White (copy formatting from row 15's A cell first)
$rec.Range("A15").Copy($rec.Range("A16"))
$rec.Cells.Item(16, 1).Value = 14
$rec.Cells.Item(16, 2).Value = "C.White"
$rec.Cells.Item(16, 3).Value = 3
$rec.Cells.Item(16, 4).Value = 2
$rec.Cells.Item(16, 5).Value = 0
$rec.Cells.Item(16, 6).Value = 0
$rec.Cells.Item(16, 7).Value = 0
$rec.Cells.Item(16, 8).Value = 0
